$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "birth year" column (old C), shifting
# old C,D,E,F,G (birth/start/end year, party, incumbent-bool) one to the right.
$ws.Columns("C").Insert()

# The freshly inserted column doesn't pick up a custom width automatically;
# give it a width close to its neighbours (best achievable via ColumnWidth).
$ws.Columns("C").ColumnWidth = 10.33

# New column C = presidency number, counting down from the 44th (Obama) to
# the 30th (Coolidge, the new row 15 added below).
$presidencyNumbers = @(44,43,42,41,40,39,38,37,36,35,34,33,32,31)
for ($i = 0; $i -lt $presidencyNumbers.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $presidencyNumbers[$i]
}

# Append row 15: Calvin Coolidge.
$ws.Range("A15").Value = "Calvin"
$ws.Range("B15").Value = "Coolidge"
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 1872
$ws.Range("E15").Value = 1923
$ws.Range("F15").Value = 1929
$ws.Range("G15").Value = "Republican"

# Old column G (now H) held an xlTRUE/xlFALSE flag; replace it with the
# "Yes"/"No" text it now represents.
for ($row = 1; $row -le 14; $row++) {
    $flag = $ws.Cells.Item($row, 8).Value2
    if ($flag) {
        $ws.Cells.Item($row, 8).Value = "Yes"
    } else {
        $ws.Cells.Item($row, 8).Value = "No"
    }
}
$ws.Range("H15").Value = "No"

$ws.Range("H13").Select() | Out-Null
